$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.654.20'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.82%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.344.41'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -2.75%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '543.81'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.22%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '136.17'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -5.46%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.25%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.522'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -9.05%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.342.79'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.68%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.104'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.76%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.156'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.51%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.28'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.04%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.338'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.32%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '24.48'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -4.18%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.767.80'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.74%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '60.543.25'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.72%  '

# Row 17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.93%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.341.23'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -2.82%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.55'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.90%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '317.61'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.84%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.07'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.08%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.54'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -3.77%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.11%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.02'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.60%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.71'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -5.44%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.30'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +7.71%  '

# Row 27
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.89'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.07%  '

# Row 28
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Bittensor'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '496.66'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -4.37%  '

# Row 29
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'Fetch.AI'
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.36'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -5.43%  '

# Row 30
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0856'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -10.10%  '

# Row 31
$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.144'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.45%  '

# Row 32
$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'PancakeSwap'
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.78'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.75%  '

# Row 33
$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'ImmutableX'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.50'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -3.95%  '

# Row 34
$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'FirstDigitalUSD'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.997'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.28%  '

# Row 35
$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.56'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -1.94%  '

# Row 36
$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'PolygonEcosystemToken'
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.375'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.12%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.49'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +2.43%  '

# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '5.20'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -6.19%  '

# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.80'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.84%  '

# Row 40
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'Monero'
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '141.14'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.51%  '

# Row 41
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.23%  '

# Row 42
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'OKB'
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '40.53'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.45%  '

# Row 43
$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '141.36'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.33%  '

# Row 44
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'Filecoin'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.54'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.16%  '

# Row 45
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'dogwifhat'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.06'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -6.64%  '

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0510'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.24%  '

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '18.93'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -8.59%  '

# Row 48
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Mantle'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.566'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -2.20%  '

# Row 49
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'Stellar'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0897'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.72%  '

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0219'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.45%  '

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '16.35'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -2.69%  '

Write-Host "Applied cryptos list update"